# Add a new data row (row 16) to the worksheet, mirroring the existing
# rows' pattern (A = index, B = shared label string, C:M = numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 16

# Copy the formatting of the cell above (bordered/bold/centered style) onto
# the new index cell, matching the pattern used by every other data row.
$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 14

$ws.Cells.Item($row, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item($row, 3).Value  = 1.150524100294847
$ws.Cells.Item($row, 4).Value  = 1.829049254762831
$ws.Cells.Item($row, 5).Value  = 0.7120580361936329
$ws.Cells.Item($row, 6).Value  = 1.150524100294847
$ws.Cells.Item($row, 7).Value  = 1.247245971645248
$ws.Cells.Item($row, 8).Value  = 0.7068087551259379
$ws.Cells.Item($row, 9).Value  = 0.8241590955535638
$ws.Cells.Item($row, 10).Value = 1.829049254762831
$ws.Cells.Item($row, 11).Value = 1.270553645478232
$ws.Cells.Item($row, 12).Value = 1.21053887288654
$ws.Cells.Item($row, 13).Value = 1.07830753559601
